$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RiskList")
$ws.Range("E6").Value = "Použití architektury (MVVM), dodržování best practises. Průběžné manuální testování. Tvorba automatických testů."
$ws.Range("G4").Value = 0.35
$ws.Range("G5").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Value = 0.55
$ws.Range("D9").Select()
